# Update the vm_pu results table (rows 2-25) for the "case with 380 kV" run.
# Column B (bus 0 / slack voltage setpoint) moves from 1.05 -> 1.02 p.u. for
# every row, and all the other bus-voltage columns (C:F and I:N) are
# refreshed with the newly recomputed load-flow results. Column G (=1) and
# column A (row index) are unchanged, and there is no column H data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$block1 = New-Object 'object[,]' 24,5
$block1[0,0] = 1.02
$block1[0,1] = 1.034680183964114
$block1[0,2] = 1.042573795731473
$block1[0,3] = 1.033833102106783
$block1[0,4] = 1.051304869849652
$block1[1,0] = 1.02
$block1[1,1] = 1.035543502118159
$block1[1,2] = 1.043268524090301
$block1[1,3] = 1.03456458517171
$block1[1,4] = 1.052231963448786
$block1[2,0] = 1.02
$block1[2,1] = 1.036102605197556
$block1[2,2] = 1.043718446497382
$block1[2,3] = 1.035038695208969
$block1[2,4] = 1.052832782471279
$block1[3,0] = 1.02
$block1[3,1] = 1.036337765276846
$block1[3,2] = 1.04390768484896
$block1[3,3] = 1.035238198928144
$block1[3,4] = 1.053085587148738
$block1[4,0] = 1.02
$block1[4,1] = 1.036377256282835
$block1[4,2] = 1.043939464088784
$block1[4,3] = 1.035271707448721
$block1[4,4] = 1.053128047039767
$block1[5,0] = 1.02
$block1[5,1] = 1.036105746975224
$block1[5,2] = 1.043720974754059
$block1[5,3] = 1.03504136025134
$block1[5,4] = 1.052836159595989
$block1[6,0] = 1.02
$block1[6,1] = 1.034971846449831
$block1[6,2] = 1.042808501398282
$block1[6,3] = 1.034080145341626
$block1[6,4] = 1.051617992423885
$block1[7,0] = 1.02
$block1[7,1] = 1.032977496001405
$block1[7,2] = 1.041203639379988
$block1[7,3] = 1.032392497164619
$block1[7,4] = 1.04947860372857
$block1[8,0] = 1.02
$block1[8,1] = 1.031650518435777
$block1[8,2] = 1.040135865731048
$block1[8,3] = 1.031271621333549
$block1[8,4] = 1.048057269925008
$block1[9,0] = 1.02
$block1[9,1] = 1.031076554071346
$block1[9,2] = 1.0396740342479
$block1[9,3] = 1.030787291967589
$block1[9,4] = 1.047443005782379
$block1[10,0] = 1.02
$block1[10,1] = 1.030863453527703
$block1[10,2] = 1.039502569439207
$block1[10,3] = 1.030607544972249
$block1[10,4] = 1.047215020097225
$block1[11,0] = 1.02
$block1[11,1] = 1.030909159979561
$block1[11,2] = 1.039539345578207
$block1[11,3] = 1.030646094294962
$block1[11,4] = 1.047263915657033
$block1[12,0] = 1.02
$block1[12,1] = 1.031058937160298
$block1[12,2] = 1.039659859272714
$block1[12,3] = 1.030772430863599
$block1[12,4] = 1.047424156745822
$block1[13,0] = 1.02
$block1[13,1] = 1.031151232504227
$block1[13,2] = 1.039734122381908
$block1[13,3] = 1.030850291515101
$block1[13,4] = 1.047522910387141
$block1[14,0] = 1.02
$block1[14,1] = 1.031688623860032
$block1[14,2] = 1.04016652709154
$block1[14,3] = 1.0313037862531
$block1[14,4] = 1.048098061675247
$block1[15,0] = 1.02
$block1[15,1] = 1.032025883741985
$block1[15,2] = 1.040437903913303
$block1[15,3] = 1.031588525046694
$block1[15,4] = 1.04845915672294
$block1[16,0] = 1.02
$block1[16,1] = 1.032222662001566
$block1[16,2] = 1.040596243636542
$block1[16,3] = 1.031754706407529
$block1[16,4] = 1.04866989136351
$block1[17,0] = 1.02
$block1[17,1] = 1.032289768484988
$block1[17,2] = 1.040650241858545
$block1[17,3] = 1.031811386535705
$block1[17,4] = 1.04874176571713
$block1[18,0] = 1.02
$block1[18,1] = 1.031989692703827
$block1[18,2] = 1.04040878254307
$block1[18,3] = 1.031557965116142
$block1[18,4] = 1.048420402837198
$block1[19,0] = 1.02
$block1[19,1] = 1.031014828869727
$block1[19,2] = 1.039624368772999
$block1[19,3] = 1.030735223614798
$block1[19,4] = 1.047376964774243
$block1[20,0] = 1.02
$block1[20,1] = 1.03040244579099
$block1[20,2] = 1.039131640163845
$block1[20,3] = 1.03021882806634
$block1[20,4] = 1.046721952442752
$block1[21,0] = 1.02
$block1[21,1] = 1.030727028816142
$block1[21,2] = 1.039392800454352
$block1[21,3] = 1.03049249375609
$block1[21,4] = 1.047069087868571
$block1[22,0] = 1.02
$block1[22,1] = 1.032006045690045
$block1[22,2] = 1.040421941081802
$block1[22,3] = 1.031571773530598
$block1[22,4] = 1.048437913697945
$block1[23,0] = 1.02
$block1[23,1] = 1.033492633120069
$block1[23,2] = 1.041618165463708
$block1[23,3] = 1.032828057509774
$block1[23,4] = 1.050030826160598

$block2 = New-Object 'object[,]' 24,6
$block2[0,0] = 1.040069724736037
$block2[0,1] = 1.039798103356771
$block2[0,2] = 1.045350100056458
$block2[0,3] = 1.036634315049783
$block2[0,4] = 1.054056741365509
$block2[0,5] = 1.041274735930047
$block2[1,0] = 1.040303883901377
$block2[1,1] = 1.040305089739655
$block2[1,2] = 1.045855991859023
$block2[1,3] = 1.037175067368833
$block2[1,4] = 1.054796160103415
$block2[1,5] = 1.041782442291747
$block2[2,0] = 1.040454282682634
$block2[2,1] = 1.040632970999041
$block2[2,2] = 1.046183047473347
$block2[2,3] = 1.037525093381449
$block2[2,4] = 1.055274911069588
$block2[2,5] = 1.04211078918013
$block2[3,0] = 1.040517241976631
$block2[3,1] = 1.040770769756521
$block2[3,2] = 1.046320470979568
$block2[3,3] = 1.037672272543342
$block2[3,4] = 1.055476247747758
$block2[3,5] = 1.042248783627651
$block2[4,0] = 1.040527797361357
$block2[4,1] = 1.040793904247259
$block2[4,2] = 1.046343540811795
$block2[4,3] = 1.037696986191238
$block2[4,4] = 1.055510057110584
$block2[4,5] = 1.04227195097202
$block2[5,0] = 1.040455125003167
$block2[5,1] = 1.040634812441311
$block2[5,2] = 1.046184884012165
$block2[5,3] = 1.037527059887267
$block2[5,4] = 1.055277601067826
$block2[5,5] = 1.042112633237459
$block2[6,0] = 1.040149090938676
$block2[6,1] = 1.03996947687439
$block2[6,2] = 1.045521127929798
$block2[6,3] = 1.036817038590146
$block2[6,4] = 1.054306568883651
$block2[6,5] = 1.041446352817718
$block2[7,0] = 1.03960129094152
$block2[7,1] = 1.038795797620812
$block2[7,2] = 1.044349335135859
$block2[7,3] = 1.035566888164178
$block2[7,4] = 1.0525978305088
$block2[7,5] = 1.040271006804977
$block2[8,0] = 1.039230400305082
$block2[8,1] = 1.038012556981017
$block2[8,2] = 1.043566752135981
$block2[8,3] = 1.034734199351462
$block2[8,4] = 1.051460335281952
$block2[8,5] = 1.039486653873634
$block2[9,0] = 1.039068459949065
$block2[9,1] = 1.037673232539397
$block2[9,2] = 1.043227570772302
$block2[9,3] = 1.034373828083464
$block2[9,4] = 1.050968200214096
$block2[9,5] = 1.039146847552387
$block2[10,0] = 1.039008107140549
$block2[10,1] = 1.03754716701253
$block2[10,2] = 1.043101537153711
$block2[10,3] = 1.034239999708129
$block2[10,4] = 1.050785462075527
$block2[10,5] = 1.039020602998012
$block2[11,0] = 1.039021062104838
$block2[11,1] = 1.037574209623615
$block2[11,2] = 1.043128573880602
$block2[11,3] = 1.034268704991385
$block2[11,4] = 1.050824657164094
$block2[11,5] = 1.039047684012706
$block2[12,0] = 1.03906347526229
$block2[12,1] = 1.037662812433606
$block2[12,2] = 1.043217153737007
$block2[12,3] = 1.034362765178034
$block2[12,4] = 1.050953093734274
$block2[12,5] = 1.03913641264885
$block2[13,0] = 1.039089580801863
$block2[13,1] = 1.037717400227733
$block2[13,2] = 1.043271724576005
$block2[13,3] = 1.034420722717411
$block2[13,4] = 1.051032236101178
$block2[13,5] = 1.039191077963905
$block2[14,0] = 1.039241119533794
$block2[14,1] = 1.038035073213617
$block2[14,2] = 1.043589255892045
$block2[14,3] = 1.034758120082783
$block2[14,4] = 1.051493005383567
$block2[14,5] = 1.039509202081867
$block2[15,0] = 1.039335816915458
$block2[15,1] = 1.038234294743742
$block2[15,2] = 1.043788350658243
$block2[15,3] = 1.034969811875098
$block2[15,4] = 1.051782144010081
$block2[15,5] = 1.039708706529417
$block2[16,0] = 1.039390922739418
$block2[16,1] = 1.038350480226951
$block2[16,2] = 1.043904448528399
$block2[16,3] = 1.035093306218628
$block2[16,4] = 1.051950832934831
$block2[16,5] = 1.039825057009338
$block2[17,0] = 1.039409690404702
$block2[17,1] = 1.038390093545292
$block2[17,2] = 1.043944029652141
$block2[17,3] = 1.035135417631125
$block2[17,4] = 1.052008358112916
$block2[17,5] = 1.039864726583134
$block2[18,0] = 1.039325670178132
$block2[18,1] = 1.038212921915877
$block2[18,2] = 1.043766992844937
$block2[18,3] = 1.034947097467805
$block2[18,4] = 1.051751118100882
$block2[18,5] = 1.039687303349684
$block2[19,0] = 1.039050991186032
$block2[19,1] = 1.037636721804311
$block2[19,2] = 1.043191070452828
$block2[19,3] = 1.034335065971507
$block2[19,4] = 1.050915270621949
$block2[19,5] = 1.039110284967868
$block2[20,0] = 1.038877127320319
$block2[20,1] = 1.037274295987325
$block2[20,2] = 1.042828696401289
$block2[20,3] = 1.033950429112683
$block2[20,4] = 1.050390103843212
$block2[20,5] = 1.038747344464652
$block2[21,0] = 1.03896940577847
$block2[21,1] = 1.037466438225122
$block2[21,2] = 1.043020822825364
$block2[21,3] = 1.034154315705656
$block2[21,4] = 1.050668469726203
$block2[21,5] = 1.038939759566465
$block2[22,0] = 1.039330255452351
$block2[22,1] = 1.038222579429661
$block2[22,2] = 1.043776643616544
$block2[22,3] = 1.034957361074942
$block2[22,4] = 1.051765137252624
$block2[22,5] = 1.039696974578246
$block2[23,0] = 1.039743916590408
$block2[23,1] = 1.039099365681724
$block2[23,2] = 1.044652521471878
$block2[23,3] = 1.035889956149193
$block2[23,4] = 1.053039293506618
$block2[23,5] = 1.040575005967355

$ws.Range("B2:F25").Value = $block1
$ws.Range("I2:N25").Value = $block2

Write-Output "applied"
